# feat: add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" sheet right after "总计", built from a copy of
#    the existing "2022-Q3" sheet (so number formats / column styling line
#    up with the other quarterly sheets), then overwrite its contents with
#    the 2022-Q4 fund-holding data.
# 2) Update the "总计" summary sheet so it gains a new first data row for
#    2022-Q4 and the remaining rows shift down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q4" worksheet
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)            # "2022-Q3" (currently 2nd sheet)
$zj = $wb.Worksheets.Item(1)            # "总计"

$q3.Copy($null, $zj)                    # paste the copy right after "总计"
$q4 = $wb.Worksheets.Item(2)            # the freshly inserted copy
$q4.Name = "2022-Q4"

# The copied sheet has 12 data rows (rows 2-13); 2022-Q4 only needs 8
# (rows 2-9), so drop the trailing 4 rows and shift the rest up.
$q4.Range("A10:H13").Delete(-4162)      # xlShiftUp

$q4Data = @(
    @(0, "015769", "天弘低碳经济混合A",           "1.12", "86.07", "4.78", "0.0535", 4),
    @(1, "015770", "天弘低碳经济混合C",           "0.99", "86.07", "4.78", "0.0473", 4),
    @(2, "015429", "中银证券专精特新股票A",         "1.18", "61.04", "2.00", "0.0236", 9),
    @(3, "011351", "金鹰年年邮益一年持有期混合A",     "3.04", "39.17", "0.73", "0.0222", 10),
    @(4, "007046", "方正富邦创新动力混合C",         "0.19", "85.11", "4.17", "0.0079", 9),
    @(5, "730001", "方正富邦创新动力混合A",         "0.17", "85.11", "4.17", "0.0071", 9),
    @(6, "015430", "中银证券专精特新股票C",         "0.14", "61.04", "2.00", "0.0028", 9),
    @(7, "011352", "金鹰年年邮益一年持有期混合C",     "0.23", "39.17", "0.73", "0.0017", 10)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = $row[0]

    # Fund code / size / position columns are stored as literal text in the
    # workbook (e.g. "015769", "1.12") even though they look numeric, so
    # force a text format before assigning to avoid Excel auto-converting
    # them (and dropping leading zeros on the fund codes).
    $c2 = $q4.Cells.Item($r, 2); $c2.NumberFormat = "@"; $c2.Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $c4 = $q4.Cells.Item($r, 4); $c4.NumberFormat = "@"; $c4.Value = $row[3]
    $c5 = $q4.Cells.Item($r, 5); $c5.NumberFormat = "@"; $c5.Value = $row[4]
    $c6 = $q4.Cells.Item($r, 6); $c6.NumberFormat = "@"; $c6.Value = $row[5]
    $c7 = $q4.Cells.Item($r, 7); $c7.NumberFormat = "@"; $c7.Value = $row[6]

    $q4.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary table
# ---------------------------------------------------------------------
# Extend the styled A column (s="2") down into the new row 5 before
# overwriting values, so the new row picks up the same formatting as its
# neighbours.
$zj.Range("A4").Copy($zj.Range("A5"))

$zj.Range("A5").Value = 3
$zj.Range("B5").Value = "2021-Q4"
$zj.Range("C5").Value = 1
$zj.Range("D5").Value = 0.42

$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2022-Q2"
$zj.Range("C4").Value = 2
$zj.Range("D4").Value = 0.02

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2022-Q3"
$zj.Range("C3").Value = 12
$zj.Range("D3").Value = 2.42

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 8
$zj.Range("D2").Value = 0.17

$zj.Activate()
